# "set up for call from python"
# - reset sheetView (drop the frozen/scrolled topLeftCell, move selection to I5)
# - rewrite row 1 to hold just two counters (A1/B1) instead of 52 zero cells
# - tweak a handful of the 0/1 "inlet" grid cells in rows 6-13 and 22
# - append a brand-new all-ones row 27
#
# NOTE: the commit also rewrites the x15ac:absPath="...\server\" ->
# "...\server\test\" attribute in xl/workbook.xml. That is Excel-internal
# "last saved from" metadata derived from the real filesystem path the file
# is saved to; it is not exposed anywhere on the Workbook/Application COM
# object model (no settable Path/FullName/AbsolutePath property reaches it),
# so it cannot be driven from this script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old scroll position and select I5 instead of Y30.
$ws.Range("I5").Select()

# Row 1: A1/B1 become counters, the rest of the row (C1:AZ1) is cleared out.
$ws.Range("C1:AZ1").ClearContents()
$ws.Range("A1").Value = 26
$ws.Range("B1").Value = 52

# Row 6: A6:H6 flip from 1 to 0.
$ws.Range("A6:H6").Value = 0

# Row 7: I7:P7 flip from 1 to 0.
$ws.Range("I7:P7").Value = 0

# Row 8: Q8 and W8:AZ8 flip from 1 to 0.
$ws.Range("Q8").Value = 0
$ws.Range("W8:AZ8").Value = 0

# Row 9: M9:N9 go from 0.5 to 1, V9 flips from 1 to 0.
$ws.Range("M9:N9").Value = 1
$ws.Range("V9").Value = 0

# Row 10: Q10 and V10 flip from 0 to 1.
$ws.Range("Q10").Value = 1
$ws.Range("V10").Value = 1

# Row 11: A11:K11 go from 0 to 1, M11:N11 go from 1 to 0.5,
# P11 flips from 0 to 1, and W11:AZ11 go from 0 to 1.
$ws.Range("A11:K11").Value = 1
$ws.Range("M11:N11").Value = 0.5
$ws.Range("P11").Value = 1
$ws.Range("W11:AZ11").Value = 1

# Row 12: L12 and O12 flip from 0 to 1.
$ws.Range("L12").Value = 1
$ws.Range("O12").Value = 1

# Row 13: M13:N13 flip from 0 to 1.
$ws.Range("M13:N13").Value = 1

# Row 22: the whole row flips from 1 to 0.
$ws.Range("A22:AZ22").Value = 0

# Row 27: brand-new row, all ones.
$ws.Range("A27:AZ27").Value = 1
